# "Products can be deleted, updated and added now. Bugfixes in adding- and
# updating-process"
#
# On the "TODO CMS" sheet:
#   - "Produkte löschen" (row 5) and "Produkte hinzufügen" (row 6) move from
#     "offen" to "done" (value + the "done" status formatting).
#   - "Löschen von Usern" (row 4) gets a follow-up remark in column D asking
#     whether it's even worth keeping.
#   - The active selection moves from C5 to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO CMS")

# Mark "Produkte löschen" and "Produkte hinzufügen" as done: copy the
# formatting of an already-"done" status cell (B2) so the same shared style
# (green "done" fill) is reused, then set the text.
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B5").Value = "done"
$ws.Range("B6").Value = "done"

# New remark on the "Löschen von Usern" row.
$ws.Range("D4").Value = "Überhaupt sinnvoll?"

# Update the sheet's active selection.
$ws.Activate() | Out-Null
$ws.Range("D5").Select() | Out-Null
